$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1117.125
$ws.Range("I70").Value = 1089.75
$ws.Range("J70").Value = 1144.5
$ws.Range("K70").Value = 3269.25
$ws.Range("L70").Value = 3433.5
$ws.Range("M70").Value = -2999.25
$ws.Range("N70").Value = -3973.5

$ws.Range("H73").Value = 1117.125
$ws.Range("I73").Value = 1089.75
$ws.Range("J73").Value = 1144.5
$ws.Range("K73").Value = 3269.25
$ws.Range("L73").Value = 3433.5
$ws.Range("M73").Value = -2333.25
$ws.Range("N73").Value = -5305.5

$ws.Range("H76").Value = 3706615.8
$ws.Range("I76").Value = 2911.4285
$ws.Range("K76").Value = 2911.4285
$ws.Range("M76").Value = -2596.4285

$ws.Range("H79").Value = 3706615.8
$ws.Range("I79").Value = 2911.4285
$ws.Range("K79").Value = 2911.4285
$ws.Range("M79").Value = -1819.4285

$ws.Range("H96").Value = 17858660
$ws.Range("J96").Value = 1827
$ws.Range("L96").Value = 5481
$ws.Range("N96").Value = -8227

$ws.Range("H99").Value = 231.55556
$ws.Range("I99").Value = 196.8
$ws.Range("K99").Value = 590.4000000000001
$ws.Range("M99").Value = 907.5999999999999

$ws.Range("H101").Value = 292
$ws.Range("I101").Value = 220
$ws.Range("J101").Value = 400
$ws.Range("K101").Value = 660
$ws.Range("L101").Value = 1200
$ws.Range("M101").Value = 962
$ws.Range("N101").Value = -4444

$ws.Range("H112").Value = 3473309
$ws.Range("J112").Value = 3704809.8
$ws.Range("L112").Value = 11114429.4
$ws.Range("N112").Value = -11116645.4

$ws.Range("H129").Value = 761.2759
$ws.Range("J129").Value = 806.9231
$ws.Range("L129").Value = 2420.7693
$ws.Range("N129").Value = -12420.7693

$ws.Range("H132").Value = 6608.2
$ws.Range("I132").Value = 7093.0835
$ws.Range("J132").Value = 4668.6665
$ws.Range("K132").Value = 21279.2505
$ws.Range("L132").Value = 14005.9995
$ws.Range("M132").Value = -18749.2505
$ws.Range("N132").Value = -19065.9995

$ws.Range("H135").Value = 14709371
$ws.Range("I135").Value = 631.29034
$ws.Range("J135").Value = 166699680
$ws.Range("K135").Value = 5681.61306
$ws.Range("L135").Value = 1500297120
$ws.Range("M135").Value = -3146.61306
$ws.Range("N135").Value = -1500302190

$ws.Range("H137").Value = 107057.84
$ws.Range("I137").Value = 126227.97
$ws.Range("K137").Value = 378683.91
$ws.Range("M137").Value = -376133.91

$ws.Range("H138").Value = 2851.9158
$ws.Range("I138").Value = 2610
$ws.Range("J138").Value = 2885.0547
$ws.Range("K138").Value = 7830
$ws.Range("L138").Value = 8655.1641
$ws.Range("M138").Value = -2690
$ws.Range("N138").Value = -18935.1641

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 2976
$ws.Range("I28").Value = 2976
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 2976
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -2784

$ws.Range("H32").Value = 7338.9126
$ws.Range("I32").Value = 5605.847
$ws.Range("J32").Value = 22936.5
$ws.Range("K32").Value = 5605.847
$ws.Range("L32").Value = 22936.5
$ws.Range("M32").Value = -5318.847
$ws.Range("N32").Value = -23510.5

$ws.Range("H45").Value = 2768.6924
$ws.Range("I45").Value = 2923
$ws.Range("J45").Value = 2614.3845
$ws.Range("K45").Value = 2923
$ws.Range("L45").Value = 2614.3845
$ws.Range("M45").Value = -2546
$ws.Range("N45").Value = -3368.3845

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0

$ws.Range("H99").Value = 2976
$ws.Range("I99").Value = 2976
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2976
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 19

$ws.Range("H122").Value = 1976.1613
$ws.Range("I122").Value = 1837.3462
$ws.Range("J122").Value = 2698
$ws.Range("K122").Value = 5512.0386
$ws.Range("L122").Value = 8094
$ws.Range("M122").Value = -3062.0386
$ws.Range("N122").Value = -12994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1895.8695
$ws.Range("I86").Value = 1700.3125
$ws.Range("J86").Value = 2342.8572
$ws.Range("K86").Value = 1700.3125
$ws.Range("L86").Value = 2342.8572
$ws.Range("M86").Value = -577.3125
$ws.Range("N86").Value = -4588.8572

$ws.Range("H89").Value = 1895.8695
$ws.Range("I89").Value = 1700.3125
$ws.Range("J89").Value = 2342.8572
$ws.Range("K89").Value = 8501.5625
$ws.Range("L89").Value = 11714.286
$ws.Range("M89").Value = -2885.5625
$ws.Range("N89").Value = -22946.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4109.622
$ws.Range("I31").Value = 1689.0741
$ws.Range("J31").Value = 7740.4443
$ws.Range("K31").Value = 1689.0741
$ws.Range("L31").Value = 7740.4443
$ws.Range("M31").Value = -1394.0741
$ws.Range("N31").Value = -8330.444299999999

$ws.Range("H34").Value = 4109.622
$ws.Range("I34").Value = 1689.0741
$ws.Range("J34").Value = 7740.4443
$ws.Range("K34").Value = 1689.0741
$ws.Range("L34").Value = 7740.4443
$ws.Range("M34").Value = -1487.0741
$ws.Range("N34").Value = -8144.4443

$ws.Range("H94").Value = 3830.8667
$ws.Range("I94").Value = 2451.5557
$ws.Range("J94").Value = 5899.8335
$ws.Range("K94").Value = 2451.5557
$ws.Range("L94").Value = 5899.8335
$ws.Range("M94").Value = -2000.5557
$ws.Range("N94").Value = -6801.8335

$ws.Range("H97").Value = 31998.5
$ws.Range("J97").Value = 31998.5
$ws.Range("L97").Value = 31998.5
$ws.Range("N97").Value = -33980.5

$ws.Range("H105").Value = 8334086
$ws.Range("I105").Value = 11364280
$ws.Range("J105").Value = 1052.5
$ws.Range("K105").Value = 11364280
$ws.Range("L105").Value = 1052.5
$ws.Range("M105").Value = -11362533
$ws.Range("N105").Value = -4546.5

$ws.Range("H109").Value = 82522650
$ws.Range("J109").Value = 82522650
$ws.Range("L109").Value = 82522650
$ws.Range("N109").Value = -82524730

$ws.Range("H141").Value = 27881.828
$ws.Range("J141").Value = 27881.828
$ws.Range("L141").Value = 27881.828
$ws.Range("N141").Value = -38241.828

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 5664.143
$ws.Range("J81").Value = 5664.143
$ws.Range("L81").Value = 16992.429
$ws.Range("N81").Value = -19238.429

$ws.Range("H84").Value = 5664.143
$ws.Range("J84").Value = 5664.143
$ws.Range("L84").Value = 50977.287
$ws.Range("N84").Value = -62209.287

$ws.Range("H131").Value = 769.47
$ws.Range("I131").Value = 345.4
$ws.Range("J131").Value = 791.7895
$ws.Range("K131").Value = 1036.2
$ws.Range("L131").Value = 2375.3685
$ws.Range("M131").Value = 4003.8
$ws.Range("N131").Value = -12455.3685

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 2875.5
$ws.Range("I41").Value = 3546.3333
$ws.Range("J41").Value = 863
$ws.Range("K41").Value = 3546.3333
$ws.Range("L41").Value = 863
$ws.Range("M41").Value = -3191.3333
$ws.Range("N41").Value = -1573

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3025
$ws.Range("I46").Value = 3033.3333
$ws.Range("K46").Value = 3033.3333
$ws.Range("M46").Value = -2845.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3711
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 4388.75
$ws.Range("K81").Value = 2000
$ws.Range("L81").Value = 8777.5
$ws.Range("M81").Value = -939
$ws.Range("N81").Value = -10899.5

$ws.Range("H84").Value = 3711
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 4388.75
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 43887.5
$ws.Range("M84").Value = -4696
$ws.Range("N84").Value = -54495.5

$ws.Range("H96").Value = 3600
$ws.Range("J96").Value = 3600
$ws.Range("L96").Value = 3600
$ws.Range("N96").Value = -6346

$ws.Range("H100").Value = 420
$ws.Range("I100").Value = 422.85715
$ws.Range("J100").Value = 416.66666
$ws.Range("K100").Value = 845.7143
$ws.Range("L100").Value = 833.33332
$ws.Range("M100").Value = -304.7143
$ws.Range("N100").Value = -1915.33332
